$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 247 (shifts existing rows 247.. down by one)
$ws.Range("A247:R247").Insert()

# Fill in the new row's data
$ws.Cells.Item(247, 1).Value = 11
$ws.Cells.Item(247, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(247, 3).Value = "Bíobío"
$ws.Cells.Item(247, 4).Value = 44588
$ws.Cells.Item(247, 5).Value = 8
$ws.Cells.Item(247, 6).Value = 100112020
$ws.Cells.Item(247, 7).Value = "Tomate"
$ws.Cells.Item(247, 8).Value = "Semiduro"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 800
$ws.Cells.Item(247, 11).Value = 5000
$ws.Cells.Item(247, 12).Value = 5500
$ws.Cells.Item(247, 13).Value = 5312
$ws.Cells.Item(247, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(247, 15).Value = "Región de Ñuble"
$ws.Cells.Item(247, 16).Value = 531
$ws.Cells.Item(247, 17).Value = 10
$ws.Cells.Item(247, 18).Value = "Hortaliza"
